# Regenerate merged AHB files
# - Rename header row labels: "_old" -> "_FV2404", "_new" -> "_FV2410"
# - Turn the data range A1:U66 into a real Excel Table (Table1)
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" headers (columns A..J = 1..10) to "_FV2404" ---
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2404")
    }
}

# --- 2. Rename the "_new" headers (columns L..U = 12..21) to "_FV2410" ---
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2410")
    }
}

# --- 3. Convert A1:U66 into an Excel Table (ListObject) named Table1 ---
$range = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $range, $false, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 4. Freeze the top (header) row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Edit applied: headers renamed, Table1 created, top row frozen."
